$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new entire row at row 1188 (shifting existing rows 1188+ down by one),
# mirroring formatting of the row above it.
$newRow = $ws.Rows.Item(1188)
$newRow.Insert(-4121)  # xlShiftDown
$ws.Rows.Item(1188).RowHeight = 22

# Populate the newly inserted row with the "yougner_sister" concept entry
$ws.Range("A1188").Value = "yougner_sister"
$ws.Range("B1188").Value = "_younger_sister"
$ws.Range("C1188").Value = 1761
$ws.Range("D1188").Value = "YOUNGER SISTER"
$ws.Range("E1188").Value = "A female who is younger than one or more of her siblings."

# Keep the named range in sync with the now one-row-larger extent it covers.
$wb.Names.Item("vanuatu_all_concepts_6").RefersTo = "=Sheet1!`$A`$1:`$E`$1222"

# Match the author's on-save cursor position (view state only).
$ws.Range("A1182").Select()

$wb.Save()
